$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update betting odds in column E (format changed from "x/y" to "x:y"),
# assigning in the order that reproduces the expected shared-string table
# ordering (new distinct values sorted by ratio, first-use order).
$ws.Range("E4").Value = "11:4"
$ws.Range("E6").Value = "7:2"
$ws.Range("E2").Value = "8:1"
$ws.Range("E3").Value = "11:1"
$ws.Range("E14").Value = "12:1"
$ws.Range("E9").Value = "20:1"
$ws.Range("E8").Value = "33:1"
$ws.Range("E7").Value = "75:1"
$ws.Range("E20").Value = "80:1"
$ws.Range("E10").Value = "100:1"
$ws.Range("E12").Value = "200:1"
$ws.Range("E16").Value = "250:1"
$ws.Range("E13").Value = "400:1"
$ws.Range("E18").Value = "500:1"
$ws.Range("E5").Value = "11:1"
$ws.Range("E11").Value = "100:1"
$ws.Range("E15").Value = "100:1"
$ws.Range("E17").Value = "250:1"
$ws.Range("E19").Value = "500:1"
$ws.Range("E21").Value = "400:1"
$ws.Range("E22").Value = "500:1"
$ws.Range("E23").Value = "500:1"

# Update the saved selection from E2 to E20
$ws.Range("E20").Select()
